$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly record per row (rows 2..607). Two new weekly
# records need to be inserted right before the old row 513, which pushes
# every row from 513 on down by two (513->515, 514->516, ..., 607->609).
$ws.Rows.Item(513).Resize(2).Insert()

# --- New row 513: Choclo / Choclero -----------------------------------
$ws.Cells.Item(513, 1).Value = 10
$ws.Cells.Item(513, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(513, 3).Value = "La Araucanía"
$ws.Cells.Item(513, 4).Value = 45015
$ws.Cells.Item(513, 5).Value = 9
$ws.Cells.Item(513, 6).Value = 100112024
$ws.Cells.Item(513, 7).Value = "Choclo"
$ws.Cells.Item(513, 8).Value = "Choclero"
$ws.Cells.Item(513, 9).Value = "Primera"
$ws.Cells.Item(513, 10).Value = 12000
$ws.Cells.Item(513, 11).Value = 500
$ws.Cells.Item(513, 12).Value = 500
$ws.Cells.Item(513, 13).Value = 500
$ws.Cells.Item(513, 14).Value = "$/unidad"
$ws.Cells.Item(513, 15).Value = "Región del Maule"
$ws.Cells.Item(513, 16).Value = 500
$ws.Cells.Item(513, 17).Value = 1
$ws.Cells.Item(513, 18).Value = "Hortaliza"

# --- New row 514: Choclo / Dulce o Americano ---------------------------
$ws.Cells.Item(514, 1).Value = 10
$ws.Cells.Item(514, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(514, 3).Value = "La Araucanía"
$ws.Cells.Item(514, 4).Value = 45015
$ws.Cells.Item(514, 5).Value = 9
$ws.Cells.Item(514, 6).Value = 100112024
$ws.Cells.Item(514, 7).Value = "Choclo"
$ws.Cells.Item(514, 8).Value = "Dulce o Americano"
$ws.Cells.Item(514, 9).Value = "Primera"
$ws.Cells.Item(514, 10).Value = 10000
$ws.Cells.Item(514, 11).Value = 280
$ws.Cells.Item(514, 12).Value = 280
$ws.Cells.Item(514, 13).Value = 280
$ws.Cells.Item(514, 14).Value = "$/unidad"
$ws.Cells.Item(514, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(514, 16).Value = 280
$ws.Cells.Item(514, 17).Value = 1
$ws.Cells.Item(514, 18).Value = "Hortaliza"
